$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect before writing values, then restore protection
$ws.Unprotect()

# Update the confidential disclaimer date from 2021-04-21 to 2021-04-22
$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-22 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) columns for each holding row (2-77)
$ws.Range("D2").Value = 0.06717372130477627
$ws.Range("E2").Value = -0.01168539325842699
$ws.Range("D3").Value = 0.04027811565919099
$ws.Range("E3").Value = -0.01575838335286528
$ws.Range("D4").Value = 0.03434018392608963
$ws.Range("E4").Value = -0.01308619234016406
$ws.Range("D5").Value = 0.0299741915495782
$ws.Range("E5").Value = 0.0005275884809847753
$ws.Range("D6").Value = 0.02729538932312056
$ws.Range("E6").Value = -0.01133715188623341
$ws.Range("D7").Value = 0.02404692231192205
$ws.Range("E7").Value = -0.02105752623887325
$ws.Range("D8").Value = 0.1749268236831203
$ws.Range("E8").Value = -0.01706484641638228
$ws.Range("D9").Value = 0.0246149098501884
$ws.Range("E9").Value = -0.008463893390959876
$ws.Range("D10").Value = 0.02295311800389171
$ws.Range("E10").Value = -0.01622214103032515
$ws.Range("D11").Value = 0.02259059315875483
$ws.Range("E11").Value = -0.007092198581560405
$ws.Range("D12").Value = 0.02047470679670009
$ws.Range("E12").Value = -0.00191141936540884
$ws.Range("D13").Value = 0.01886439044339542
$ws.Range("E13").Value = -0.009297520661157077
$ws.Range("D14").Value = 0.01721250233818641
$ws.Range("E14").Value = -0.008280377431157282
$ws.Range("D15").Value = 0.01742850772344147
$ws.Range("E15").Value = -0.001484780994803425
$ws.Range("D16").Value = 0.01591690930549462
$ws.Range("E16").Value = -0.02211874272409786
$ws.Range("D17").Value = 0.0146607315649542
$ws.Range("E17").Value = -0.01083569405099161
$ws.Range("D18").Value = 0.01467390993011119
$ws.Range("E18").Value = -0.006191318546740932
$ws.Range("D19").Value = 0.01324293914626122
$ws.Range("E19").Value = -0.01641954423325709
$ws.Range("D20").Value = 0.01252344034312933
$ws.Range("E20").Value = -0.01303571428571415
$ws.Range("D21").Value = 0.01190401724631083
$ws.Range("E21").Value = 0.04151444702756568
$ws.Range("D22").Value = 0.0128317342189232
$ws.Range("E22").Value = -0.006768953068592043
$ws.Range("D23").Value = 0.01223910713125726
$ws.Range("E23").Value = 0.0006786739754633953
$ws.Range("D24").Value = 0.0127410031715393
$ws.Range("E24").Value = -0.0056919335021689
$ws.Range("D25").Value = 0.0117149676079678
$ws.Range("E25").Value = -0.01377511888326421
$ws.Range("D26").Value = 0.0100064524670236
$ws.Range("E26").Value = -0.05341378925019558
$ws.Range("D27").Value = 0.01024318382657102
$ws.Range("E27").Value = -0.04165302144249516
$ws.Range("D28").Value = 0.01071241349503967
$ws.Range("E28").Value = -0.01427027027027017
$ws.Range("D29").Value = 0.010493532830114
$ws.Range("E29").Value = 0.0005822604645142615
$ws.Range("D30").Value = 0.01011379624139328
$ws.Range("E30").Value = -0.003158809128958384
$ws.Range("D31").Value = 0.009920912896822762
$ws.Range("E31").Value = -0.01773940345368918
$ws.Range("D32").Value = 0.010040037330954
$ws.Range("E32").Value = -0.004665629860031162
$ws.Range("D33").Value = 0.009519332169493336
$ws.Range("E33").Value = -0.01479188166494672
$ws.Range("D34").Value = 0.009159442997387844
$ws.Range("E34").Value = -0.003112982970151945
$ws.Range("D35").Value = 0.009057929651239137
$ws.Range("E35").Value = -0.003659289304294044
$ws.Range("D36").Value = 0.008812732190439347
$ws.Range("E36").Value = 0.007340946166394913
$ws.Range("D37").Value = 0.008476683878936051
$ws.Range("E37").Value = -0.0007773302240122737
$ws.Range("D38").Value = 0.008914804618746231
$ws.Range("E38").Value = -0.0328307262269526
$ws.Range("D39").Value = 0.008840207085069683
$ws.Range("E39").Value = -0.02251454591449542
$ws.Range("D40").Value = 0.007794683514841785
$ws.Range("E40").Value = -0.01732697362016944
$ws.Range("D41").Value = 0.007344582443070705
$ws.Range("E41").Value = -0.01017855977728965
$ws.Range("D42").Value = 0.007614139912190994
$ws.Range("E42").Value = -0.02364343931272483
$ws.Range("D43").Value = 0.007959572817063661
$ws.Range("E43").Value = -0.02334985650926158
$ws.Range("D44").Value = 0.007262477234698658
$ws.Range("E44").Value = -0.005542725173210195
$ws.Range("D45").Value = 0.007673961703115764
$ws.Range("E45").Value = 0.007909910284964905
$ws.Range("D46").Value = 0.007161363232948647
$ws.Range("E46").Value = -0.0133832976445396
$ws.Range("D47").Value = 0.007593613610097982
$ws.Range("E47").Value = -0.0160292818376877
$ws.Range("D48").Value = 0.007101781048663095
$ws.Range("E48").Value = -0.008794619762027889
$ws.Range("D49").Value = 0.006992919765578367
$ws.Range("E49").Value = -0.007024156244646251
$ws.Range("D50").Value = 0.006738537383608554
$ws.Range("E50").Value = -0.006779661016949157
$ws.Range("D51").Value = 0.006498211924472864
$ws.Range("E51").Value = -0.01535748085692168
$ws.Range("D52").Value = 0.00648159919748708
$ws.Range("E52").Value = -0.01997461584907512
$ws.Range("D53").Value = 0.005475011706132112
$ws.Range("E53").Value = -0.02297592997811815
$ws.Range("D54").Value = 0.006100385034491185
$ws.Range("E54").Value = -0.02042419481539659
$ws.Range("D55").Value = 0.005588824859760679
$ws.Range("E55").Value = 0.01736334405144691
$ws.Range("D56").Value = 0.005688435329849598
$ws.Range("E56").Value = -0.006416189901529923
$ws.Range("D57").Value = 0.005704874342021948
$ws.Range("E57").Value = -0.005936047488379947
$ws.Range("D58").Value = 0.005647049273090662
$ws.Range("E58").Value = -0.01731160896130346
$ws.Range("D59").Value = 0.005078462718226264
$ws.Range("E59").Value = -0.01188959660297251
$ws.Range("D60").Value = 0.005006820333100071
$ws.Range("E60").Value = -0.003190403266972952
$ws.Range("D61").Value = 0.004749482802579901
$ws.Range("E61").Value = -0.05660377358490576
$ws.Range("D62").Value = 0.004566503199097061
$ws.Range("E62").Value = -0.03130738959335366
$ws.Range("D63").Value = 0.004519620166690066
$ws.Range("E63").Value = -0.008058245564430577
$ws.Range("D64").Value = 0.004237523283450698
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0.004123470523182913
$ws.Range("E65").Value = -0.003912605562872873
$ws.Range("D66").Value = 0.003813579269794254
$ws.Range("E66").Value = -0.006785624528776069
$ws.Range("D67").Value = 0.003977230604380173
$ws.Range("E67").Value = -0.008795710584974969
$ws.Range("D68").Value = 0.003333647171440528
$ws.Range("E68").Value = -0.006181269316466453
$ws.Range("D69").Value = 0.003549772360015197
$ws.Range("E69").Value = 0.0004724940938238831
$ws.Range("D70").Value = 0.003048595139650806
$ws.Range("E70").Value = -0.01498559077809802
$ws.Range("D71").Value = 0.003173470133123272
$ws.Range("E71").Value = -0.009890898108648849
$ws.Range("D72").Value = 0.002453651854472423
$ws.Range("E72").Value = -0.03321831971615508
$ws.Range("D73").Value = 0.002032263644967638
$ws.Range("E73").Value = -0.0002358027117311812
$ws.Range("D74").Value = 0.002046001092282806
$ws.Range("E74").Value = -0.009856735761408442
$ws.Range("D75").Value = 0.001527971538293413
$ws.Range("E75").Value = -0.0153677277716795
$ws.Range("D76").Value = 0.001405612414532736
$ws.Range("E76").Value = -0.0100573896244105
$ws.Range("D77").Value = 0.9999999999999998
$ws.Range("E77").Value = -0.01189562901721641

# Restore sheet protection
$ws.Protect()
